# Commit message: "Added units to description texts in excel"
#
# The cost-calculation sheet had two variable-machine-cost line-item
# labels ("Treibstoffe (l/Std.)" and "Schmierstoffkosten") repeated for
# each of the three machine blocks (Grundmaschine rows 51-52, Maschine 1
# rows 57-58, Maschine 2 rows 63-64). The labels are updated to state the
# unit consistently as "(l/h)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsverfahren")

# Grundmaschine (base machine) variable cost block
$ws.Range("B51").Value = "Treibstoffe (l/h)"
$ws.Range("B52").Value = "Schmierstoffkosten (l/h)"

# Maschine 1 variable cost block
$ws.Range("B57").Value = "Treibstoffe (l/h)"
$ws.Range("B58").Value = "Schmierstoffkosten (l/h)"

# Maschine 2 variable cost block
$ws.Range("B63").Value = "Treibstoffe (l/h)"
$ws.Range("B64").Value = "Schmierstoffkosten (l/h)"

# Restore the active selection recorded at save time.
$ws.Range("B65").Select() | Out-Null
